$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "![](datasets/ETH/seq_eth/reference.png)"
$ws.Range("B2").Value = "[ETH](datasets/ETH)"
$ws.Range("A3").Value = "![](datasets/UCY/zara01/reference.png)"
$ws.Range("B3").Value = "[UCY](datasets/UCY)"
$ws.Range("A4").Value = "![](datasets/PETS-2009/reference.jpg)"
$ws.Range("B4").Value = "[PETS 2009](datasets/PETS-2009)"
$ws.Range("A5").Value = "![](datasets/SDD/coupa/video3/reference.jpg)"
$ws.Range("B5").Value = "[SDD](datasets/SDD)"
$ws.Range("A6").Value = "![](datasets/GC/reference.jpg)"
$ws.Range("B6").Value = "[GC](datasets/GC)"
$ws.Range("A7").Value = "![](datasets/HERMES/reference.png)"
$ws.Range("B7").Value = "[HERMES](datasets/HERMES)"
$ws.Range("A8").Value = "![](datasets/Waymo/reference.jpg)"
$ws.Range("B8").Value = "[Waymo](datasets/Waymo)"
$ws.Range("A9").Value = "![](datasets/KITTI/reference.jpg)"
$ws.Range("B9").Value = "[KITTI](datasets/KITTI)"
$ws.Range("A10").Value = "![](datasets/InD/reference.png)"
$ws.Range("B10").Value = "[inD](datasets/InD)"
$ws.Range("A11").Value = "![](datasets/L-CAS/reference.png)"
$ws.Range("B11").Value = "[L-CAS](datasets/L-CAS)"
$ws.Range("A12").Value = "![](datasets/VIRAT/reference.png)"
$ws.Range("B12").Value = "[VIRAT](datasets/VIRAT)"
$ws.Range("A13").Value = "![](datasets/VRU/reference.png)"
$ws.Range("B13").Value = "[VRU](datasets/VRU)"
$ws.Range("A14").Value = "![](datasets/Edinburgh/reference.jpg)"
$ws.Range("B14").Value = "[Edinburgh](datasets/Edinburgh)"
$ws.Range("A15").Value = "![](datasets/Town-Center/reference.jpg)"
$ws.Range("B15").Value = "[Town Center](datasets/Town-Center)"
$ws.Range("A16").Value = "![](datasets/ATC/reference.png)"
$ws.Range("B16").Value = "[ATC](datasets/ATC)"
$ws.Range("A17").Value = "![](datasets/City-Scapes/reference.png)"
$ws.Range("B17").Value = "[City Scapes](datasets/City-Scapes)"
$ws.Range("A18").Value = "![](datasets/Forking-Paths-Garden/reference.png)"
$ws.Range("B18").Value = "[Forking Paths Garden](datasets/Forking-Paths-Garden)"
$ws.Range("A19").Value = "![](datasets/NuScenes/reference.png)"
$ws.Range("B19").Value = "[nuScenes](datasets/NuScenes)"
$ws.Range("A20").Value = "![](datasets/Argoverse/reference.jpg)"
$ws.Range("B20").Value = "[Argoverse](datasets/Argoverse)"
$ws.Range("A21").Value = "![](datasets/Wild-Track/reference.jpg)"
$ws.Range("B21").Value = "[Wild Track](datasets/Wild-Track)"
$ws.Range("A22").Value = "![](datasets/DUT/reference.png)"
$ws.Range("B22").Value = "[DUT](datasets/DUT)"
$ws.Range("A23").Value = "![](datasets/CITR/reference.png)"
$ws.Range("B23").Value = "[CITR](datasets/CITR)"
$ws.Range("A24").Value = "![](datasets/Ko-PER/reference.png)"
$ws.Range("B24").Value = "[Ko-PER](datasets/Ko-PER)"
$ws.Range("A25").Value = "![](datasets/TRAF/reference.png)"
$ws.Range("B25").Value = "[TRAF](datasets/TRAF)"
$ws.Range("A26").Value = "![](datasets/ETH-Person/reference.png)"
$ws.Range("B26").Value = "[ETH-Person](datasets/ETH-Person)"

$ws.Range("C26").Value = "Multi-Person Data Collected from Mobile Platforms"
